$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

function Set-CellText($row, $col, $oldText, $newText) {
    $cell = $tbl.Cell($row, $col)
    $cellRng = $cell.Range
    # Build a fresh Range over the cell's start/end so the Find is confined
    # to this cell only (Wrap = wdFindStop, Replace = wdReplaceOne).
    $rng = $d.Range($cellRng.Start, $cellRng.End)
    $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1)
}

$map = @(
    @{ Row = 1;  Col = 1; Old = "20÷7=2, 6";  New = "24÷7=3, 3" },
    @{ Row = 1;  Col = 2; Old = "83÷5=16, 3"; New = "48÷3=16, 0" },
    @{ Row = 1;  Col = 3; Old = "74÷6=12, 2"; New = "97÷9=10, 7" },
    @{ Row = 1;  Col = 4; Old = "67÷5=13, 2"; New = "96÷6=16, 0" },
    @{ Row = 1;  Col = 5; Old = "78÷6=13, 0"; New = "34÷6=5, 4" },

    @{ Row = 5;  Col = 1; Old = "33÷2=16, 1"; New = "28÷9=3, 1" },
    @{ Row = 5;  Col = 2; Old = "34÷8=4, 2";  New = "21÷6=3, 3" },
    @{ Row = 5;  Col = 3; Old = "45÷8=5, 5";  New = "40÷4=10, 0" },
    @{ Row = 5;  Col = 4; Old = "34÷8=4, 2";  New = "77÷6=12, 5" },
    @{ Row = 5;  Col = 5; Old = "95÷9=10, 5"; New = "92÷4=23, 0" },

    @{ Row = 9;  Col = 1; Old = "39÷4=9, 3";  New = "48÷6=8, 0" },
    @{ Row = 9;  Col = 2; Old = "17÷3=5, 2";  New = "67÷3=22, 1" },
    @{ Row = 9;  Col = 3; Old = "35÷5=7, 0";  New = "80÷6=13, 2" },
    @{ Row = 9;  Col = 4; Old = "20÷4=5, 0";  New = "71÷3=23, 2" },
    @{ Row = 9;  Col = 5; Old = "48÷9=5, 3";  New = "71÷5=14, 1" },

    @{ Row = 13; Col = 1; Old = "95÷4=23, 3"; New = "29÷2=14, 1" },
    @{ Row = 13; Col = 2; Old = "18÷3=6, 0";  New = "18÷2=9, 0" },
    @{ Row = 13; Col = 3; Old = "53÷3=17, 2"; New = "57÷6=9, 3" },
    @{ Row = 13; Col = 4; Old = "18÷4=4, 2";  New = "90÷5=18, 0" },
    @{ Row = 13; Col = 5; Old = "60÷5=12, 0"; New = "64÷7=9, 1" },

    @{ Row = 17; Col = 1; Old = "80÷9=8, 8";  New = "44÷5=8, 4" },
    @{ Row = 17; Col = 2; Old = "25÷3=8, 1";  New = "38÷2=19, 0" },
    @{ Row = 17; Col = 3; Old = "86÷4=21, 2"; New = "47÷7=6, 5" },
    @{ Row = 17; Col = 4; Old = "92÷5=18, 2"; New = "26÷5=5, 1" }
)

foreach ($item in $map) {
    Set-CellText $item.Row $item.Col $item.Old $item.New
}
